$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.178.11'
$ws.Range("E2").Value = '  +2.36%  '

$ws.Range("D3").Value = '2.943.47'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.995'
$ws.Range("E4").Value = '  -0.48%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.09'
$ws.Range("E5").Value = '  +1.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.98'
$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").Value = '2.929.06'
$ws.Range("E8").Value = '  +0.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  -0.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  +5.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.142'
$ws.Range("E11").Value = '  +1.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("E13").Value = '  +0.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.24'
$ws.Range("E14").Value = '  -0.09%  '

$ws.Range("E15").Value = '  +0.03%  '

$ws.Range("D16").Value = '3.411.93'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").Value = '61.646.60'
$ws.Range("E17").Value = '  +1.58%  '

$ws.Range("D18").Value = '2.941.38'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.68'
$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '439.17'
$ws.Range("E20").Value = '  +1.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.47'
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.673'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.07'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.19'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.86'
$ws.Range("E25").Value = '  +0.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").Value = '  -1.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.84'
$ws.Range("E27").Value = '  +0.79%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  -6.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.91'
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.40'
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  -2.37%  '

$ws.Range("D35").Value = '0.0₃0868'
$ws.Range("E35").Value = '  +1.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.988'
$ws.Range("E36").Value = '  -1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.60'
$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.42'
$ws.Range("E38").Value = '  -0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.98'
$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.87'
$ws.Range("E40").Value = '  -3.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.52'
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.118'
$ws.Range("E42").Value = '  -1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.277'
$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.97'
$ws.Range("E44").Value = '  -5.45%  '

$ws.Range("D45").Value = '2.689.73'
$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '360.73'
$ws.Range("E48").Value = '  -3.60%  '

$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.89'
$ws.Range("E50").Value = '  -3.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.104'
$ws.Range("E51").Value = '  -1.73%  '
